# Page Object Model Framework
# Restructure the TestCase_Scenarios workbook: add new sheets for the
# expanded test scenarios (Comman verification, Products Details,
# Order History, Address), drop & recreate the blank "Sheet4", reorder
# all tabs, and populate every sheet's S.No/TC-No/Description(/Comment)
# table.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------
# 1. Drop the old blank "Sheet4" and add the new sheets (appended at the
#    end, one at a time, so each gets a fresh, predictable sheetId).
# ---------------------------------------------------------------------
$oldSheet4 = $wb.Worksheets.Item("Sheet4")
$oldSheet4.Delete()

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$orderHistory = $wb.Worksheets.Add($null, $last)
$orderHistory.Name = "Order History"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$address = $wb.Worksheets.Add($null, $last)
$address.Name = "Address"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$commanVerification = $wb.Worksheets.Add($null, $last)
$commanVerification.Name = "Comman verification"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$productsDetails = $wb.Worksheets.Add($null, $last)
$productsDetails.Name = "Products Details"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet4 = $wb.Worksheets.Add($null, $last)
$newSheet4.Name = "Sheet4"

# ---------------------------------------------------------------------
# 2. Reorder tabs into the final layout:
#    Comman verification, Home, Login, MyAccount, Products Details,
#    Sheet4, Order History, Address
# ---------------------------------------------------------------------
$commanVerification = $wb.Worksheets.Item("Comman verification")
$home = $wb.Worksheets.Item("Home")
$commanVerification.Move($home)

$productsDetails = $wb.Worksheets.Item("Products Details")
$myAccount = $wb.Worksheets.Item("MyAccount")
$productsDetails.Move($null, $myAccount)

$newSheet4 = $wb.Worksheets.Item("Sheet4")
$productsDetails = $wb.Worksheets.Item("Products Details")
$newSheet4.Move($null, $productsDetails)

# ---------------------------------------------------------------------
# 3. Populate "Comman verification"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Comman verification")
$ws.Range("A1").Value = "S.No"
$ws.Range("B1").Value = "TC-No"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Comment"
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("C2").Value = "TC001"
$ws.Range("C3").Value = "Check Women, Dresses, Tshirt Header is displaying on each page."

# ---------------------------------------------------------------------
# 4. Populate "Home" (headers unchanged, just rewritten)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Home")
$ws.Range("A1").Value = "S.No"
$ws.Range("B1").Value = "TC-No"
$ws.Range("C1").Value = "Description"
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("C2").Value = "Verify Logo"
$ws.Range("C3").Value = "Login Page verification"
$ws.Range("C4").Value = "Check Women, Dresses, Tshirt Header is displaying on each page."

# ---------------------------------------------------------------------
# 5. Populate "Login"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Login")
$ws.Range("A1").Value = "S.No"
$ws.Range("B1").Value = "TC-No"
$ws.Range("C1").Value = "Description"
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "TC001"

# ---------------------------------------------------------------------
# 6. Populate "MyAccount"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("MyAccount")
$ws.Range("A1").Value = "S.No"
$ws.Range("B1").Value = "TC-No"
$ws.Range("C1").Value = "Description"
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A2").Font.Bold = $true
$ws.Range("B2").Font.Bold = $true
$ws.Range("C2").Value = "Verify Comman verification points"
$ws.Range("C2").Font.Bold = $true
$ws.Range("C3").Value = "Check if Cart is empty after login."
$ws.Range("C4").Value = "Search Sleevs & Verify result."
$ws.Range("C5").Value = "Verify all labels in MyAccount Section"
$ws.Range("C6").ClearContents()
# Column widths (11 / 60, both bestFit) are already correct on this
# pre-existing sheet - leave them untouched.

# ---------------------------------------------------------------------
# 7. Populate "Products Details"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Products Details")
$ws.Range("A1").Value = "S.No"
$ws.Range("B1").Value = "TC-No"
$ws.Range("C1").Value = "Description"
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("C2").Value = "Verifu Mouse hover on Women / Dresses / T-Shirts"
$ws.Range("C3").Value = "verify Tooltip on mousehover"
$ws.Range("C4").Value = "Verify price on mouse hover on a item"
$ws.Range("C5").Value = "Add Item to Cart & Verify count in Cart"
$ws.Range("C6").Value = "MouseHover on cart & verify item in cart"
$ws.Range("C7").Value = "Verify total in cart after addign multiple items"
$ws.Range("C8").Value = "Verify end to end order process while not logegd in"
$ws.Range("C9").Value = "Verify end to end order process while not logegdin"
$ws.Columns.Item(3).ColumnWidth = 46.16666666666667

# ---------------------------------------------------------------------
# 8. "Sheet4" stays blank (freshly created, nothing to populate).
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 9. Populate "Order History"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Order History")
$ws.Range("A1").Value = "S.No"
$ws.Range("B1").Value = "TC-No"
$ws.Range("C1").Value = "Description"
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("C2").Value = "Verify Comman verification points"
$ws.Range("C2").Font.Bold = $true
$ws.Range("C3").Value = "Verify Order History for user with no orders"
$ws.Range("C4").Value = "Craete a order & verify history."

# ---------------------------------------------------------------------
# 10. Populate "Address"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Address")
$ws.Range("A1").Value = "S.No"
$ws.Range("B1").Value = "TC-No"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Comment"
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A2").Font.Bold = $true
$ws.Range("B2").Font.Bold = $true
$ws.Range("C2").Value = "Verify Comman verification points"
$ws.Range("C2").Font.Bold = $true
$ws.Range("D2").Font.Bold = $true
$ws.Range("C3").Value = "Verify Address is same as entered."
$ws.Range("D3").Value = "Can be done using JDBC connector of read from file"
$ws.Range("C4").Value = "Verify user can update address"
$ws.Range("C5").Value = "Verify user can delete address"
$ws.Range("C6").Value = "Add a new address & verify its displayed."
$ws.Range("C7").Value = 'Go back to Account using "Back to your account"'
$ws.Range("C8").Value = 'Go back to Home page using "Home"'
$ws.Columns.Item(3).ColumnWidth = 31.30625

# ---------------------------------------------------------------------
# 11. Selections per sheet (mirrors the diff's per-sheet cursor state).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("MyAccount").Range("C25").Select()
$wb.Worksheets.Item("Login").Range("D15").Select()
$wb.Worksheets.Item("Home").Range("G20").Select()
$wb.Worksheets.Item("Comman verification").Range("C7").Select()
$wb.Worksheets.Item("Order History").Range("G16").Select()
$wb.Worksheets.Item("Address").Range("L16").Select()
$wb.Worksheets.Item("Sheet4").Range("A1").Select()

# "Products Details" is the active tab in the final workbook (activeTab=4).
$wb.Worksheets.Item("Products Details").Range("G24").Select()
